# DW_P5_Plan_tests_acceptation.xlsx
# "plan d'acceptation, message de suppression d'article"
#
# Row 6 (test case 5, "Sur la page cart, l'utilisateur peut changer la
# quantité d'un produit.") — the expected-result cell (E6) gains a sentence
# about the on-screen quantity being clamped to 100.
#
# Row 7 (test case 6, "Supprimer un produit du panier.") — the expected-
# result cell (E7) is replaced with a longer message describing the new
# confirmation dialog behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the expected-result text for test case 5 (row 6) ---
$ws.Cells.Item(6, 5).Value = 'OK / Si l''utilisateur cherche à outrepasser la limite fixée, un message d''alerte le rappelera à l''ordre ! La quantité de l''article dans le Localstorage ne pourra pas franchir le plafond de 100. Sur la page, la quantité sera initialisée à 100.'

# --- Update the expected-result text for test case 6 (row 7) ---
$ws.Cells.Item(7, 5).Value = 'OK / Une message de confirmation est générée pour valider la suppression de l''article. En cas de décrémentation de l''Input jusqu''à 0, si l''utilisateur ne valide pas la suppression, la quantité du produit est réinitialisée à 1.'

# --- Row heights grew to fit the longer wrapped text ---
$ws.Rows.Item(6).RowHeight = 151.2
$ws.Rows.Item(7).RowHeight = 129.6

# --- The visible top-left cell of the scrolled view moved from D6 to A7 ---
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 8
